$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    2 = 10.83;
    3 = 8.92;
    4 = 12.35;
    5 = 15.29;
    6 = 3.11;
    7 = 8.44;
    8 = 7.64;
    9 = 5.63;
    10 = 24.72;
    11 = 1.21;
    12 = 3.23;
    13 = 1.27;
    14 = 6.74;
    15 = 3.87;
    16 = 4.77;
    17 = 16.13;
    18 = 3.14;
    19 = 12.93;
    21 = 1.36;
    22 = 7.3;
    23 = 2.96;
    24 = 3.32;
    26 = 5.22;
    28 = 4.62;
    29 = 0.78;
    30 = 1.37;
    31 = 2.26;
    32 = 1.41;
    33 = 2.78;
    37 = 5.47;
    38 = 26.63;
    39 = 2.12;
    40 = 0.93;
    41 = 2.23;
    42 = 3.73;
    43 = 3.37;
    44 = 0.96;
    46 = 3.35;
    47 = 8.01;
    48 = 4.32;
    51 = 6.07;
    53 = 3.69;
    54 = 4.02;
    55 = 1.35;
    56 = 11.1;
    57 = 8.35;
    58 = 17.04;
    59 = 1.57;
    60 = 9.42;
    61 = 8.539999999999999;
    63 = 1.88;
    64 = 5.55;
    66 = 1.08;
    67 = 0.11
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = $priceUpdates[$row]
}
